# Stand-up Meeting Sprint2.xlsx
#
# Harish Chowdary Bala's block (rows 20-22 on Sheet1) was still the blank
# "1) / 2) / 3)" placeholder text, formatted with the plain thin-bordered
# grey style (same as the still-empty Manoj Kumar Gude block above it).
# This edit fills in Harish's actual standup notes and gives his block the
# same "filled in" look already used for the other completed members
# (Saibabu Devarapalli, rows 14-16): a copy of that block's cell formatting
# (fill/border/font/alignment), auto-sized row heights for the new text,
# and the real write-up text in place of the placeholders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the "filled in" box formatting from Saibabu's block (A14:D16)
# onto Harish's block (A20:D22) -------------------------------------------
$ws.Range("A14:D16").Copy()
$ws.Range("A20").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 20: "1) What did I accomplish?" -----------------------------------
$ws.Range("B20").Value = "Present I don't accomplish anything just discuss about project with team members."
$ws.Range("C20").Value = ".I learn about the firebase database it helps me to do design these project"
$ws.Range("D20").Value = "I learn about MongoDB and Firebase Databse"

# --- Row 21: "2) What will I do today?" ------------------------------------
$ws.Range("B21").Value = "Today we discuss about the data of everyone is collected and what type of technologies are used in these design."
$ws.Range("C21").Value = "Yesterday we discuss about the design of these project and we are working on these assigned work."
$ws.Range("D21").Value = "2)I will continue the which work assigned by team discussion."

# --- Row 22: "3) What obstacles are impeding my progress?" -----------------
$ws.Range("B22").Value = "We are confusing about team roles and work structure."
$ws.Range("C22").Value = "I don't know about andriod studio, So I am working on that."
$ws.Range("D22").Value = "I have no experience and I learning it will take some time."

# --- Resize the newly filled rows to fit their (now much longer) text ------
$ws.Rows(20).RowHeight = 47
$ws.Rows(21).RowHeight = 62.5
$ws.Rows(22).RowHeight = 31.5

# --- Leave the view where the author left it when saving -------------------
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("E17").Select()
